$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "js in browser" execution time value: 1h 2m -> 1h 19m
$ws.Range("C2").Value = "1h 19m"

# Column C's cells (C1:C5) previously used the numeric-format cell style
# (same one used by column B/D). Re-format them to match column A's
# style (General number format, same border/font/alignment) by copying
# column A's formatting onto C1:C5.
$ws.Range("A1:A5").Copy()
$ws.Range("C1:C5").PasteSpecial(-4122)
